$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45192.45601851852
$ws.Range("B3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").Value = 45192.45664351852
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 'Hossam.Ibrahim'
$ws.Range("E3").Value = 'Hossam Tabana'
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = 'Geospatial Maps'
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = 'Sometimes'
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = 'Yes, simple calculations'
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = ""
$ws.Range("R3").Value = 'Beginner level'
$ws.Range("S3").Value = ""
$ws.Range("T3").Value = ""
$ws.Range("U3").Value = 'Basic transformations only'
$ws.Range("V3").Value = ""
$ws.Range("W3").Value = ""
$ws.Range("X3").Value = 'Yes, basic automation'
$ws.Range("Y3").Value = ""
$ws.Range("Z3").Value = ""
$ws.Range("AA3").Value = 'Occasionally'
$ws.Range("AB3").Value = ""
$ws.Range("AC3").Value = ""
$ws.Range("AD3").Value = 'Very Important'
$ws.Range("AE3").Value = ""
$ws.Range("AF3").Value = ""
$ws.Range("AG3").Value = 'Yes, as static files'
$ws.Range("AH3").Value = ""
$ws.Range("AI3").Value = ""
$ws.Range("AJ3").Value = 'Occasionally'
$ws.Range("AK3").Value = ""
$ws.Range("AL3").Value = ""
$ws.Range("AM3").Value = 'Yes, it''s essential'
$ws.Range("AN3").Value = ""
$ws.Range("AO3").Value = ""
$ws.Range("AP3").Value = 'Not concerned'
$ws.Range("AQ3").Value = ""
$ws.Range("AR3").Value = ""
$ws.Range("AS3").Value = 'Just exploring'
$ws.Range("AT3").Value = ""
$ws.Range("AU3").Value = ""
$ws.Range("AV3").Value = 'Yes, basic trend lines'
$ws.Range("AW3").Value = ""
$ws.Range("AX3").Value = ""
$ws.Range("AY3").Value = 'Beginner'
$ws.Range("AZ3").Value = ""
$ws.Range("BA3").Value = ""
$ws.Range("BB3").Value = 'Definitely'
$ws.Range("BC3").Value = ""
$ws.Range("BD3").Value = ""
$ws.Range("BE3").Value = 'I know what it is but haven''t used it'
$ws.Range("BF3").Value = ""
$ws.Range("BG3").Value = ""
$ws.Range("BH3").Value = 'Possibly'
$ws.Range("BI3").Value = ""
$ws.Range("BJ3").Value = ""
$ws.Range("BK3").Value = 'Yes, to multiple formats'
$ws.Range("BL3").Value = ""
$ws.Range("BM3").Value = ""
$ws.Range("BN3").Value = 'Using Power BI workspaces'
$ws.Range("BO3").Value = ""
$ws.Range("BP3").Value = ""
$ws.Range("BQ3").Value = 'Basic Training'
